$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
Write-Host "ActiveSheet: $($ws.Name)"
foreach ($sheet in $wb.Worksheets) {
    Write-Host $sheet.Name
}
